$wb = $excel.ActiveWorkbook

# The "OldImportLogic" worksheet holds the TestAutomation_1 -> TestAutomation_2 testdata update
$ws = $wb.Worksheets.Item("OldImportLogic")

$ws.Range("A2").Value = "OldImportLogic_2 - TestAutomation_2"
$ws.Range("B2").Value = "OldImportLogic_2 - TestAutomation_2_radio_button"
$ws.Range("H2").Value = "StandardExcelReport-OldImportLogic_2-TestAutomation_2-Quality of Life-2023_"
$ws.Range("H3").Value = "ExcelReport-OldImportLogic_2-TestAutomation_2-Quality of Life-"
$ws.Range("H4").Value = "WordReport-OldImportLogic_2 - TestAutomation_2-Quality of Life-"

# Update the selected range on that sheet to reflect the new active selection
$ws.Activate()
$ws.Range("H2:H4").Select()
